# Implements "full duplex messages": adds First Name / Last Name columns
# (with sample data) to the "Data" sheet, next to the existing Primary Key
# column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Header row (row 1): C1/D1 get the same look as B1 ("Primary Key") ---
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$ws.Range("C1").Value = "First Name"
$ws.Range("D1").Value = "Last Name"

# Propagate that same text-cell formatting down through the data rows
# (C2:D7) before filling in the names.
$ws.Range("C1:D1").Copy()
$ws.Range("C2:D7").PasteSpecial(-4122)

# --- Data rows 2-7: Primary Key numbers in column B, names in C/D ---
$primaryKeys = @(14, 25, 109, 120, 123, 124)
$firstNames  = @("John", "Jane", "Michael", "Emily", "David", "Olivia")
$lastNames   = @("Doe", "Smith", "Johnson", "Brown", "Lee", "Garcia")

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.NumberFormat = "General"
    $bCell.Value = $primaryKeys[$i]

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = $firstNames[$i]

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $lastNames[$i]
}
